$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# "CasesTab" row is renamed to "ParticipantsTab" (rest of that row's
# content - query / accession file names - stays the same).
$ws.Range("A2").Value = "ParticipantsTab"

# Move the active selection from B3 to A2.
$ws.Range("A2").Select()
